$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (row 47) ---
$ws.Range("B47").Value = 4
$ws.Range("C47").Value = 850
$ws.Range("D47").Value = 7032
$ws.Range("E47").Formula = "=C47/D47"

# --- Update sheet view (selection) ---
$ws.Range("H42").Select() | Out-Null

# --- Update chart source ranges to include the new row (and some headroom) ---
$chartObjs = $ws.ChartObjects()

$co1 = $chartObjs.Item(1)
$chart1 = $co1.Chart
$ser1 = $chart1.SeriesCollection(1)
$ser1.Formula = "=SERIES(Sheet1!`$B`$1,Sheet1!`$A`$2:`$A`$60,Sheet1!`$B`$2:`$B`$60,1)"

$co2 = $chartObjs.Item(2)
$chart2 = $co2.Chart
$ser2 = $chart2.SeriesCollection(1)
$ser2.Formula = "=SERIES(Sheet1!`$E`$1,Sheet1!`$A`$2:`$A`$60,Sheet1!`$E`$2:`$E`$60,1)"
